$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A2").Value = "Version: $newVersion"
$aboutSheet.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Zhaozhuang Coal Mine No. 2 Well, China, M2299, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 7; $r++) {
    $dataSheet.Cells.Item($r, 19).Value = $newVersion
}
